$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column cells whose new text happens to look like a plain number need
# their format forced to Text first, otherwise Excel would coerce the value to
# a numeric cell and drop meaningful trailing zeros (e.g. "1.00" -> 1).

$ws.Range("D2").Value = "70.113.52"
$ws.Range("E2").Value = "  +0.66%  "
$ws.Range("D3").Value = "3.797.36"
$ws.Range("E3").Value = "  +4.78%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "617.63"
$ws.Range("E5").Value = "  +4.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.33"
$ws.Range("E6").Value = "  -4.36%  "
$ws.Range("D7").Value = "3.797.22"
$ws.Range("E7").Value = "  +4.89%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.537"
$ws.Range("E9").Value = "  +0.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.169"
$ws.Range("E10").Value = "  +3.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.34"
$ws.Range("E11").Value = "  -3.03%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.495"
$ws.Range("E12").Value = "  -0.69%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "41.04"
$ws.Range("E13").Value = "  +4.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000257"
$ws.Range("E14").Value = "  +0.12%  "
$ws.Range("D15").Value = "4.415.38"
$ws.Range("E15").Value = "  +4.36%  "
$ws.Range("D16").Value = "3.785.35"
$ws.Range("E16").Value = "  +4.09%  "
$ws.Range("D17").Value = "70.131.21"
$ws.Range("E17").Value = "  +0.57%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.61"
$ws.Range("E19").Value = "  +0.74%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "515.57"
$ws.Range("E20").Value = "  +1.34%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.68"
$ws.Range("E21").Value = "  -3.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.61"
$ws.Range("E22").Value = "  +3.79%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.730"
$ws.Range("E23").Value = "  -2.54%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.53"
$ws.Range("E24").Value = "  +4.60%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "88.15"
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.35"
$ws.Range("E26").Value = "  -1.60%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.11"
$ws.Range("E27").Value = "  +2.41%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000138"
$ws.Range("E28").Value = "  +28.48%  "
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.50"
$ws.Range("E30").Value = "  -1.82%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.86"
$ws.Range("E31").Value = "  -5.04%  "
$ws.Range("E32").Value = "  +3.22%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.88"
$ws.Range("E33").Value = "  -0.75%  "
$ws.Range("E34").Value = "  -2.07%  "
$ws.Range("E35").Value = "  -0.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.23"
$ws.Range("E36").Value = "  +0.89%  "
$ws.Range("E37").Value = "  +3.64%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.341"
$ws.Range("E38").Value = "  +1.53%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.18"
$ws.Range("E39").Value = "  +3.62%  "
$ws.Range("E40").Value = "  +3.25%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "51.37"
$ws.Range("E41").Value = "  +1.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "44.55"
$ws.Range("E42").Value = "  -5.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.81"
$ws.Range("E43").Value = "  -0.21%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "425.29"
$ws.Range("E44").Value = "  +5.26%  "
$ws.Range("D45").Value = "3.069.79"
$ws.Range("E45").Value = "  -2.35%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.77"
$ws.Range("E46").Value = "  -0.88%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0366"
$ws.Range("E47").Value = "  -0.03%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "27.79"
$ws.Range("E48").Value = "  -0.44%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "135.90"
$ws.Range("E49").Value = "  -0.51%  "
$ws.Range("B50").Value = "USDe"
$ws.Range("C50").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  -0.03%  "
$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.50"
$ws.Range("E51").Value = "  +1.46%  "
